# Automatische test-sync: 2025-08-03 15:02:50
$wb = $excel.ActiveWorkbook

# --- Logs sheet: append testmail #14 as a new row (row 22) -----------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A22").Value = "Heb je de CE-certificaten van dit product?"
$logs.Range("B22").Value = "mailmind.test@zohomail.eu"
$logs.Range("C22").Value = "Testmail #14: Heb je de CE-certificaten van dit product?"
$logs.Range("D22").Value = "Overig"
$logs.Range("E22").Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Range("F22").Value = "2025-08-03 15:02:02"
$logs.Range("G22").Value = "Ja"
$logs.Range("H22").Value = "Ja"
$logs.Range("I22").Value = "Nee"
$logs.Range("J22").Value = "Nee"

# --- Extend the conditional-formatting ranges to cover the new row ---------
$logs.Range("D2:D21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D22"))
$logs.Range("G2:G21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G22"))
$logs.Range("H2:H21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H22"))
$logs.Range("I2:I21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I22"))
$logs.Range("J2:J21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J22"))

# --- Dashboard sheet: "Overig" now has 7 hits and outranks "Intern verzoek" -
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A2").Value = "Overig"
$dash.Range("B2").Value = 7
$dash.Range("A3").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B3").Value = 6
